$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WBS_Week04")
$ws.Activate()

$ws.Range("S1").Value = "Estimated Cost (RM)"

$costs = @(180,120,180,240,180,180,240,240,150,150,100,100,100,150,100,100)
for ($i = 0; $i -lt $costs.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 19).Value = $costs[$i]
}

$ws.Range("S17").Select()

